# Slide 8 ("Conclusions") has a text box whose 4th bullet currently reads:
#   "Ammenities vs Occupancy: top ten ammenities in Berlin are: ..."
# Per the diff it should become:
#   "Ammenities vs Occupancy: top eight amenities in Berlin are: ..."
# i.e. the bold "top ten ammenities" (originally 2 runs) turns into
# "top eight amenities" split across 4 runs: "top ", "eight", " ", "amenities".

$p = $ppt.ActivePresentation

# Locate the slide/shape/paragraph that contains the target phrase instead of
# hard-coding indices, so the script is resilient to re-ordering.
$targetShape = $null
$targetSlide = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $t = $shape.TextFrame.TextRange.Text
            if ($t -like "*top ten ammenities*") {
                $targetShape = $shape
                $targetSlide = $slide
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange
$paraCount = $tr.Paragraphs().Count
$para = $null
for ($pi = 1; $pi -le $paraCount; $pi++) {
    $candidate = $tr.Paragraphs($pi, 1)
    if ($candidate.Text -like "*top ten ammenities*") {
        $para = $candidate
    }
}

# Step 1: "top ten ammenities" -> replace "ten" with "eight".
# Editing this as a sub-range (instead of rewriting the whole run) makes the
# host split the original "top ten " run into "top " / "eight" / " ", each
# inheriting the original run's formatting (bold, lang, fonts).
$full = $para.Text
$idx = $full.IndexOf("top ten ammenities")
$tenRange = $para.Characters($idx + 5, 3)
$tenRange.Text = "eight"

# Step 2: fix the "ammenities" typo -> "amenities" (same run, text only, so
# its existing formatting - including the spell-check err flag - is kept).
$full = $para.Text
$idx = $full.IndexOf("ammenities")
$wordRange = $para.Characters($idx + 1, 10)
$wordRange.Text = "amenities"
